# Depreciation Calculator - fill in the two missing "Asset Cost" totals and
# the Year-on-Year schedule's year numbers / running-balance seed so the
# Straight Line / Diminishing Balance sections (and the depreciation
# schedule table) actually compute instead of showing blank/"" results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Depreciation Calculator")

# --- Straight Line Method block: Asset Cost = Asset Price + Additional Asset Cost
$ws.Range("D8").Formula = "=D6+D7"

# --- Diminishing Balance Method block: same total, mirrored lower on the sheet
$ws.Range("D20").Formula = "=D18+D19"

# --- Depreciation schedule: seed the running book value from the Asset Cost
$ws.Range("D26").Formula = "=D8"

# --- Depreciation schedule: number the remaining years (1 and 2/3 were already
# filled in; continue the sequence through year 10)
$ws.Range("B29").Value = 4
$ws.Range("B30").Value = 5
$ws.Range("B31").Value = 6
$ws.Range("B32").Value = 7
$ws.Range("B33").Value = 8
$ws.Range("B34").Value = 9
$ws.Range("B35").Value = 10
